$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 60,3
$data[0,0] = 2104
$data[0,1] = 2515
$data[0,2] = 2611
$data[1,0] = 886.08300000000008
$data[1,1] = 1925.7290000000005
$data[1,2] = 2129.4700000000003
$data[2,0] = 11543.705000000004
$data[2,1] = 6488.0490000000018
$data[2,2] = 6737.002000000004
$data[3,0] = 6659.7100000000019
$data[3,1] = 3682.746000000001
$data[3,2] = 4209.5439999999999
$data[4,0] = 4494.0480000000007
$data[4,1] = 2247.5009999999997
$data[4,2] = 2209.3029999999999
$data[5,0] = 0
$data[5,1] = 0
$data[5,2] = 0
$data[6,0] = 3489
$data[6,1] = 2788
$data[6,2] = 2807
$data[7,0] = 769.35000000000014
$data[7,1] = 2061.8949999999991
$data[7,2] = 2227.2549999999992
$data[8,0] = 10928.412999999993
$data[8,1] = 6145.0419999999995
$data[8,2] = 6425.9440000000004
$data[9,0] = 6190.8819999999969
$data[9,1] = 3433.9150000000004
$data[9,2] = 4066.0949999999993
$data[10,0] = 4529.7630000000008
$data[10,1] = 2226.753999999999
$data[10,2] = 2279.1679999999992
$data[11,0] = 0
$data[11,1] = 0
$data[11,2] = 0
$data[12,0] = 3586
$data[12,1] = 2166
$data[12,2] = 1995
$data[13,0] = 1680.8300000000006
$data[13,1] = 2365.848
$data[13,2] = 2592.3199999999993
$data[14,0] = 10898.198999999999
$data[14,1] = 6496.1960000000036
$data[14,2] = 6838.2979999999998
$data[15,0] = 6418.5089999999991
$data[15,1] = 3538.2750000000015
$data[15,2] = 4137.1259999999975
$data[16,0] = 4330.2969999999987
$data[16,1] = 2277.2980000000002
$data[16,2] = 2393.0169999999994
$data[17,0] = 0
$data[17,1] = 0
$data[17,2] = 0
$data[18,0] = 1373.362499999999
$data[18,1] = 974.58799999999962
$data[18,2] = 904.50499999999965
$data[19,0] = 2226.7299999999991
$data[19,1] = 1618.7199999999998
$data[19,2] = 1595.3999999999992
$data[20,0] = 11257.601999999997
$data[20,1] = 5283.1869999999981
$data[20,2] = 5273.3289999999988
$data[21,0] = 6446.7100000000028
$data[21,1] = 2889.52
$data[21,2] = 3212.8600000000006
$data[22,0] = 4436.0550000000012
$data[22,1] = 1788.8079999999998
$data[22,2] = 1895.6379999999997
$data[23,0] = 0
$data[23,1] = 0
$data[23,2] = 0
$data[24,0] = 1437.6405000000002
$data[24,1] = 1305.6855000000005
$data[24,2] = 1296.9269999999997
$data[25,0] = 2168.4419999999991
$data[25,1] = 1919.04
$data[25,2] = 1950.3590000000004
$data[26,0] = 9275.2019999999975
$data[26,1] = 5731.1030000000019
$data[26,2] = 5661.070999999999
$data[27,0] = 5620.2959999999975
$data[27,1] = 3200.3049999999994
$data[27,2] = 3448.8949999999995
$data[28,0] = 3701.6649999999995
$data[28,1] = 1943.5560000000007
$data[28,2] = 1983.7219999999991
$data[29,0] = 0
$data[29,1] = 0
$data[29,2] = 0
$data[30,0] = 1512.3250000000003
$data[30,1] = 1711.5965000000003
$data[30,2] = 1658.5349999999992
$data[31,0] = 1903.3499999999999
$data[31,1] = 1415.1059999999993
$data[31,2] = 1425.0379999999998
$data[32,0] = 16518.094000000001
$data[32,1] = 10803.777
$data[32,2] = 10427.597999999998
$data[33,0] = 5041.8729999999978
$data[33,1] = 3122.7190000000005
$data[33,2] = 3308.4050000000007
$data[34,0] = 4795.8240000000023
$data[34,1] = 2852.4750000000004
$data[34,2] = 2852.4190000000003
$data[35,0] = 5.8399999999996908
$data[35,1] = 133.69000000000005
$data[35,2] = 121.76000000000067
$data[36,0] = 5059.6164999999946
$data[36,1] = 3905.7699999999977
$data[36,2] = 3782.1229999999991
$data[37,0] = 1969.172
$data[37,1] = 892.43800000000022
$data[37,2] = 945.06500000000028
$data[38,0] = 26598.50299999999
$data[38,1] = 14668.986999999999
$data[38,2] = 14325.897999999994
$data[39,0] = 5868.4619999999932
$data[39,1] = 3155.5130000000022
$data[39,2] = 3470.2359999999999
$data[40,0] = 8294.1510000000035
$data[40,1] = 4466.217999999998
$data[40,2] = 4292.9000000000015
$data[41,0] = 0
$data[41,1] = 0
$data[41,2] = 0
$data[42,0] = 5346.5099999999984
$data[42,1] = 4057.3584999999989
$data[42,2] = 4199.5489999999991
$data[43,0] = 1541.1400000000006
$data[43,1] = 774.31000000000063
$data[43,2] = 781.94999999999982
$data[44,0] = 23306.753000000026
$data[44,1] = 14927.280000000004
$data[44,2] = 14929.637999999992
$data[45,0] = 5499.7400000000034
$data[45,1] = 3429.6760000000004
$data[45,2] = 3680.5840000000017
$data[46,0] = 7589.8439999999946
$data[46,1] = 4775.7479999999987
$data[46,2] = 4709.6369999999997
$data[47,0] = 0
$data[47,1] = 0
$data[47,2] = 0
$data[48,0] = 3887.4975000000031
$data[48,1] = 2863.0859999999971
$data[48,2] = 3005.5939999999991
$data[49,0] = 2235.9100000000017
$data[49,1] = 1762.1320000000005
$data[49,2] = 1902.9730000000004
$data[50,0] = 16759.309999999994
$data[50,1] = 9248.1730000000025
$data[50,2] = 9299.146999999999
$data[51,0] = 4869.9800000000023
$data[51,1] = 2704.0450000000005
$data[51,2] = 2941.7550000000001
$data[52,0] = 5089.1450000000059
$data[52,1] = 2495.0159999999996
$data[52,2] = 2607.3530000000005
$data[53,0] = 0
$data[53,1] = 0
$data[53,2] = 0
$data[54,0] = 3381.9580000000005
$data[54,1] = 2380.5339999999983
$data[54,2] = 2428.5495000000019
$data[55,0] = 2284.5360000000001
$data[55,1] = 1891.0830000000001
$data[55,2] = 2068.4569999999999
$data[56,0] = 12252.739999999998
$data[56,1] = 6782.1049999999932
$data[56,2] = 6758.3839999999964
$data[57,0] = 4328.6750000000002
$data[57,1] = 2233.5149999999994
$data[57,2] = 2605.3999999999987
$data[58,0] = 3628.4829999999993
$data[58,1] = 1841.2869999999996
$data[58,2] = 1956.3139999999999
$data[59,0] = 0
$data[59,1] = 0
$data[59,2] = 0

$ws.Range("A1:C60").Value = $data

$ws.Columns.Item(1).ColumnWidth = 9.7109375
$ws.Columns.Item(2).ColumnWidth = 9.7109375
$ws.Columns.Item(3).ColumnWidth = 9.7109375

Write-Host "done"